# Trade #84 closed at 2026-02-17 08:59:31 - unknown UNKNOWN +0.000%
#
# Updates the "Summary", "Strategy Status", "All Trades" and "MarketMaking"
# sheets of the live-trading-results workbook to reflect the newly closed
# trade #84 (MarketMaking / DOWN, CLOSED, early_exit).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.53   # Current Capital
$summary.Range("B4").Value = 0.54      # Total P&L $
$summary.Range("B6").Value = 84        # Total Trades
$summary.Range("B7").Value = 37        # Winning Trades
$summary.Range("B9").Value = 44.05     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet (MarketMaking row, row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.53     # Capital
$status.Range("D4").Value = 84         # Trades
$status.Range("E4").Value = 0.54       # P&L $
$status.Range("F4").Value = 0.53       # P&L %
$status.Range("G4").Value = 44.05      # Win Rate %

# ---------------------------------------------------------------------
# 3. Append the newly closed trade as row 85 on both the "All Trades"
#    and "MarketMaking" sheets (identical row on each).
# ---------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(85, 1).Value = 84                 # Trade #

    # Date / Time columns look like dates/times to Excel's auto-detection,
    # so force them to remain literal text by pre-setting the number
    # format to Text before writing the value.
    $ws.Cells.Item(85, 2).NumberFormat = "@"
    $ws.Cells.Item(85, 2).Value = "2026-02-17"        # Date
    $ws.Cells.Item(85, 3).Value = "08:59:25"          # Time

    $ws.Cells.Item(85, 4).Value = "MarketMaking"      # Strategy
    $ws.Cells.Item(85, 5).Value = "DOWN"              # Side
    $ws.Cells.Item(85, 6).Value = 0.95                # Entry Price
    $ws.Cells.Item(85, 7).Value = 0.96                # Exit Price
    $ws.Cells.Item(85, 8).Value = "CLOSED"            # Status
    $ws.Cells.Item(85, 9).Value = 1.0526              # P&L %
    $ws.Cells.Item(85, 10).Value = 0.01               # P&L $
    $ws.Cells.Item(85, 11).Value = 100.53             # Capital After
    $ws.Cells.Item(85, 12).Value = 0                  # Entry Slippage (bps)
    $ws.Cells.Item(85, 13).Value = 0                  # Exit Slippage (bps)
    $ws.Cells.Item(85, 14).Value = 0.6                # Confidence
    $ws.Cells.Item(85, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item(85, 16).Value = "early_exit"       # Exit Reason
    $ws.Cells.Item(85, 17).Value = 0.14               # Duration (min)
}
